$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): rows 2,4,5,7 -> column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 375
$ws1.Range("F4").Value = 294
$ws1.Range("F5").Value = 4299
$ws1.Range("F7").Value = 460

# Sheet "全部类型" (All types): rows 2,4,5,9 -> column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 375
$ws4.Range("F4").Value = 294
$ws4.Range("F5").Value = 4299
$ws4.Range("F9").Value = 460
